# The workbook currently has three sheets:
#   strategy_id-0, strategy_id-5004, strategy_id-5008
#
# Target state (per the diff):
#   strategy_id-0, strategy_id-5004, strategy_id-5007, strategy_id-5009
#
# i.e. the old "strategy_id-5008" sheet is duplicated to a brand-new
# trailing sheet named "strategy_id-5009" (identical data/formatting),
# and the original sheet is renamed "strategy_id-5007" (same sheetId,
# same position, same contents).

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("strategy_id-5008")

# Duplicate the sheet, placing the copy after the last existing sheet.
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# The copy is now the last sheet in the workbook - name it.
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "strategy_id-5009"

# Rename the original sheet.
$src.Name = "strategy_id-5007"
